$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = 1234
$ws.Range("A3").Value = "new data"

$ws.Range("C1").Select()
